# TC_CreateOpportunity.xlsx - "this is my third commit"
#
# Adds a third data row (row 3) to Sheet1 that duplicates row 2 (the
# "apptestmbob4@netapp.com" / "Apple Inc." / "Amit Jain" record) but with a
# distinct Opportunity Name ("Test for Automation-Second"), including a
# fresh mailto: hyperlink on the new row's "User Name" cell (A3), mirroring
# the hyperlink already present on A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create the hyperlink on A3 first (while the cell is still blank/default
# styled) so that the subsequent row copy below re-applies the exact same
# look (border + Hyperlink style) row 2 already uses for its own link cell.
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:apptestmbob4@netapp.com")

# Duplicate row 2 into row 3 (values, shared-string reuse, styles, border).
$ws.Range("A2:D2").Copy($ws.Range("A3:D3"))

# This new row records a second/different opportunity for the same
# contact, account and owner.
$ws.Range("C3").Value = "Test for Automation-Second"
